# LOB1245.xlsx edit: remove the old "Docentes responsaveis" value-only row,
# which re-aligns every following row (and its original, already-correct
# row height) up by one, then patch the handful of B/C cells whose text
# actually changed content in this revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the stray row that only carried "7455355 - Robson da Silva Rocha"
#    in B13/C13 (no label in A13). Deleting it shifts rows 14-22 up to
#    13-21, which lines every remaining label + its custom row height back
#    into the positions the new layout expects.
$ws.Rows.Item(13).Delete()

# 2) Cascade the now-misaligned long-answer cells down one slot at a time,
#    working from the bottom up so each source is read before it gets
#    overwritten. Use Copy so the shared-string text stays marked as text
#    (not re-inferred as a date/number) and the destination keeps its own
#    column style untouched.
$ws.Range("B20").Copy($ws.Range("B21"))
$ws.Range("C20").Copy($ws.Range("C21"))

$ws.Range("B19").Copy($ws.Range("B20"))
$ws.Range("C19").Copy($ws.Range("C20"))

$ws.Range("B18").Copy($ws.Range("B19"))
$ws.Range("C18").Copy($ws.Range("C19"))

# 3) New content for "Objetivos:" (row 10) -- replace the long paragraph
#    with the teacher id/name string.
$ws.Range("B10").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C10").Value = "7455355 - Robson da Silva Rocha"

# 4) "Docentes responsaveis:" (row 18) now takes the same teacher string;
#    copy from B10/C10 so it's written as text with B18/C18's own style.
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))

# 5) "Programa resumido:" (row 13) becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 6) "Programa:" (row 15) becomes the activation date "01/01/2012" -- copy
#    from B8/C8 (already "01/01/2012" as text) to avoid it being parsed as
#    a date value/format.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
